$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tan order incorrect calculation - update values to include units / degree symbol
# Order matches the order values were entered so shared-string table indices line up
$ws.Range("E5").Value = "36.87°"
$ws.Range("G5").Value = "36.87°"

$ws.Range("E3").Value = "47°"

$ws.Range("C4").Value = "3m"
$ws.Range("D4").Value = "4m"
$ws.Range("G4").Value = "5m"

$ws.Range("D5").Value = "4mm"
$ws.Range("C5").Value = "5mm"

# Update selection to match
$ws.Range("I6").Select()
